$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (cosinor result for series 1)
$ws.Range("E2").Value2 = 24.26000000000035
$ws.Range("H2").Value2 = [double]"1.531342102931251e-16"
$ws.Range("K2").Value2 = 50.22414245239935
$ws.Range("L2").Value2 = "[44.959762059781056, 55.48852284501764]"
$ws.Range("O2").Value2 = 1.478026573760964
$ws.Range("P2").Value2 = "[1.3648160276856558, 1.591237119836272]"
$ws.Range("S2").Value2 = 58.45073679990821
$ws.Range("T2").Value2 = "[55.03996062791778, 61.86151297189865]"
$ws.Range("W2").Value2 = 18.55319319319347
$ws.Range("X2").Value2 = 18.11607607607634
$ws.Range("Y2").Value2 = 18.99031031031059

# Row 3 (cosinor result for series 2)
$ws.Range("E3").Value2 = 25.15000000000049
$ws.Range("H3").Value2 = [double]"1.531342102931251e-16"
$ws.Range("K3").Value2 = 55.33006606541628
$ws.Range("L3").Value2 = "[49.92949632555175, 60.73063580528081]"
$ws.Range("O3").Value2 = -3.107000542289005
$ws.Range("P3").Value2 = "[-3.2076321388003888, -3.0063689457776204]"
$ws.Range("S3").Value2 = 58.44647473511552
$ws.Range("T3").Value2 = "[54.76153271128146, 62.131416758949584]"
$ws.Range("W3").Value2 = 12.43653653653678
$ws.Range("X3").Value2 = 12.03373373373397
$ws.Range("Y3").Value2 = 12.83933933933959
